$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: total_registros updated (name unchanged)
$ws.Range("B3").Value = 119

# Rows 5-14 are re-sorted (descending by total_registros) with updated counts.
$ws.Range("A5").Value = "ROMERO CHANAME YOSSELY TRINIDAD"
$ws.Range("B5").Value = 102

$ws.Range("A6").Value = "VALLE MAGALLAN EDUAR"
$ws.Range("B6").Value = 102

$ws.Range("A7").Value = "ZAVALETA MANAY JORGE LUIS"
$ws.Range("B7").Value = 100

$ws.Range("A8").Value = "HIDALGO CUBAS LUISA YVONE"
$ws.Range("B8").Value = 97

$ws.Range("A9").Value = "CAMACHO LINARES JUDITH ARLETT"
$ws.Range("B9").Value = 93

$ws.Range("A10").Value = "SEVERINO AVALOS MARJORIE ISABEL"
$ws.Range("B10").Value = 91

$ws.Range("A11").Value = "HUMPIRE CASTILLO IRWIN DEIMER"
$ws.Range("B11").Value = 91

$ws.Range("A12").Value = "SENADOR ARBOLEDA GIANCARLOS EXEBIO"
$ws.Range("B12").Value = 91

$ws.Range("A13").Value = "BALLENA ESQUÉN ASTRID CAROLINA"
$ws.Range("B13").Value = 89

$ws.Range("B14").Value = 85
